$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line entries (line7, line8) are being added to the table.
# Shift the existing "extr*" rows (8-15) down two rows to (10-17),
# carrying their values/format along, which makes room for the new
# rows at 8-9.
$ws.Range("A8:E15").Copy()
$ws.Range("A10:E17").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Rows 16/17 are brand new (beyond the sheet's previous extent) and the
# bulk paste above doesn't carry the bordered/bold index-column style
# that far, so copy that formatting explicitly from the row above.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill the new row 8 with the line7 entry.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $false

# Fill the new row 9 with the line8 entry.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Renumber the sequential index column for the shifted rows (10-17).
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# The in_service flag changed for what are now extr1/extr2 (rows 10/11).
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(11, 5).Value = $true
